{"js": "// Remove the trailing \"Ver no Jupiter ...\" line, the copyright/footer line\n// that follows it, and the blank paragraph that separates them from the\n// preceding \"LOQ4073: ...\" requirement line.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph (\"LOQ4073: ...\") so the removal is resilient\n// to any unrelated structural changes elsewhere in the document.\nconst anchorIndex = items.findIndex(\n  (p) => p.text.trim() === \"LOQ4073: Qu\u00edmica Geral II (Requisito fraco)\"\n);\n\nif (anchorIndex === -1) {\n  throw new Error('Anchor paragraph \"LOQ4073: ...\" not found.');\n}\n\n// The paragraphs immediately following the anchor are expected to be:\n//   1) an empty paragraph\n//   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3) the \"\u00a9 2020 ...\" copyright/footer paragraph\n// Collect exactly those (by matching expected text) and delete them.\nconst toDelete = [];\nlet idx = anchorIndex + 1;\n\nif (items[idx] && items[idx].text.trim() === \"\") {\n  toDelete.push(items[idx]);\n  idx++;\n}\n\nif (items[idx] && items[idx].text.trim() === \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n  toDelete.push(items[idx]);\n  idx++;\n}\n\nif (\n  items[idx] &&\n  items[idx].text.trim() ===\n    \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n) {\n  toDelete.push(items[idx]);\n  idx++;\n}\n\n// Delete from the end backwards so earlier deletions don't invalidate\n// references to the paragraphs that still need removing.\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" line, the copyright/footer line\n# that follows it, and the blank paragraph that separates them from the\n# preceding \"LOQ4073: ...\" requirement line.\n\n$d = $word.ActiveDocument\n\n$anchorText = \"LOQ4073: Qu\u00edmica Geral II (Requisito fraco)\"\n$targets = @(\n    \"\",\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n)\n\n# Locate the anchor paragraph so the removal is resilient to any unrelated\n# structural changes elsewhere in the document.\n$anchorIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t -eq $anchorText) {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Anchor paragraph 'LOQ4073: ...' not found.\"\n}\n\n# Collect the ranges that exactly match the expected run of paragraphs\n# following the anchor (empty line, \"Ver no Jupiter ...\", \"\u00a9 2020 ...\").\n$rangesToDelete = @()\n$idx = $anchorIndex + 1\nforeach ($target in $targets) {\n    if ($idx -gt $d.Paragraphs.Count) { break }\n    $p = $d.Paragraphs.Item($idx)\n    $t = $p.Range.Text.Trim()\n    if ($t -eq $target) {\n        $rangesToDelete += , $p.Range\n        $idx++\n    } else {\n        break\n    }\n}\n\n# Delete from the end backwards so earlier deletions don't invalidate\n# references to the paragraphs that still need removing.\nfor ($i = $rangesToDelete.Count - 1; $i -ge 0; $i--) {\n    $rangesToDelete[$i].Delete()\n}\n"}
